$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "1.00", "0.630") that must
# stay text, exactly as in the source inlineStr cells. Force text format before
# assigning so Excel does not silently coerce them to numbers.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.361.17'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.929.71'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '486.30'
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.84'
$ws.Range('E6').Value = '  +3.34%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +1.22%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.732'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('E10').Value = '  +3.72%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000356'
$ws.Range('E11').Value = '  +6.08%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.04'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('E13').Value = '  +3.11%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.554.16'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.72'
$ws.Range('E15').Value = '  -2.30%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.918.59'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.03'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('E19').Value = '  -2.17%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '68.443.66'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '442.83'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('E22').Value = '  +3.48%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.91'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.68'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.42'
$ws.Range('E25').Value = '  +13.91%  '
$ws.Range('E26').Value = '  +15.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.62'
$ws.Range('E27').Value = '  +1.91%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '38.85'
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.88'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '714.09'
$ws.Range('E30').Value = '  -1.77%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.63'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.130'
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.90'
$ws.Range('E33').Value = '  +3.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0₃0924'
$ws.Range('E34').Value = '  +15.52%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '41.97'
$ws.Range('E35').Value = '  -2.80%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.09'
$ws.Range('E36').Value = '  +12.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '60.50'
$ws.Range('E37').Value = '  +5.42%  '
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.401'
$ws.Range('E39').Value = '  +19.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  +13.34%  '
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.20'
$ws.Range('E43').Value = '  +4.28%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.95'
$ws.Range('E44').Value = '  +6.06%  '
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0₆0367'
$ws.Range('E46').Value = '  +49.91%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.42'
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '145.75'
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.16'
$ws.Range('E51').Value = '  -1.23%  '
